# Reorder slides: move the "Methodology" slide from position 2 to the
# end of the deck (it becomes the last slide).
$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text.Trim() -eq "Methodology") {
                $targetIndex = $i
            }
            break
        }
    }
    if ($targetIndex -ne -1) { break }
}

if ($targetIndex -eq -1) {
    # Fallback: the "Methodology" slide is originally the 2nd slide.
    $targetIndex = 2
}

$s = $p.Slides.Item($targetIndex)
$s.MoveTo($p.Slides.Count)
